$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the "municipio"/"Casos"/"Óbitos" header row) and shift
# everything below it up by one row.
$ws.Rows.Item(2).Delete()
